# Update "想去人数" (interested-count) figures across sheets to match the
# latest scrape output (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 175
$wsExpo.Range("F3").Value = 490
$wsExpo.Range("F4").Value = 20
$wsExpo.Range("F9").Value = 1369

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 89
$wsShow.Range("F3").Value = 39

# Sheet "全部类型" (All types - combined view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 175
$wsAll.Range("F3").Value = 89
$wsAll.Range("F4").Value = 490
$wsAll.Range("F5").Value = 20
$wsAll.Range("F10").Value = 1369
$wsAll.Range("F11").Value = 39
